$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 491
$ws1.Range("F4").Value = 60
$ws1.Range("F5").Value = 50
$ws1.Range("F6").Value = 5
$ws1.Range("F7").Value = 265
$ws1.Range("F9").Value = 1063
$ws1.Range("F10").Value = 14974
$ws1.Range("F11").Value = 187
$ws1.Range("F14").Value = 5944
$ws1.Range("F15").Value = 607
$ws1.Range("F17").Value = 53
$ws1.Range("F18").Value = 91
$ws1.Range("F19").Value = 1247
$ws1.Range("F20").Value = 22
$ws1.Range("F22").Value = 201
$ws1.Range("F23").Value = 822
$ws1.Range("F24").Value = 2959
$ws1.Range("F25").Value = 100
$ws1.Range("F26").Value = 10771
$ws1.Range("F28").Value = 84
$ws1.Range("F29").Value = 126

# Sheet 4: 全部类型 (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 491
$ws4.Range("F5").Value = 60
$ws4.Range("F6").Value = 50
$ws4.Range("F7").Value = 5
$ws4.Range("F8").Value = 265
$ws4.Range("F10").Value = 1063
$ws4.Range("F11").Value = 14975
$ws4.Range("F12").Value = 187
$ws4.Range("F15").Value = 5944
$ws4.Range("F16").Value = 607
$ws4.Range("F18").Value = 53
$ws4.Range("F19").Value = 91
$ws4.Range("F20").Value = 1247
$ws4.Range("F21").Value = 22
$ws4.Range("F23").Value = 201
$ws4.Range("F24").Value = 822
$ws4.Range("F25").Value = 2959
$ws4.Range("F26").Value = 100
$ws4.Range("F28").Value = 10771
$ws4.Range("F30").Value = 84
$ws4.Range("F31").Value = 126
